$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.64"
$ws.Range("E2").Value = "'5.60%"
$ws.Range("D3").Value = "'31.73"
$ws.Range("E3").Value = "'7.35%"
$ws.Range("D4").Value = "'5.209"
$ws.Range("E4").Value = "'2.71%"
$ws.Range("D5").Value = "'0.07375"
$ws.Range("E5").Value = "'9.53%"
$ws.Range("D6").Value = "'7.845"
$ws.Range("E6").Value = "'6.83%"
$ws.Range("D7").Value = "'3.731"
$ws.Range("E7").Value = "'8.27%"
$ws.Range("D8").Value = "'1.504"
$ws.Range("E8").Value = "'8.20%"
$ws.Range("D9").Value = "'0.9079"
$ws.Range("E9").Value = "'-0.77%"
$ws.Range("D10").Value = "'0.01668"
$ws.Range("E10").Value = "'2,475.85%"
$ws.Range("D11").Value = "'0.1678"
$ws.Range("E11").Value = "'5.51%"
$ws.Range("D12").Value = "'0.07455"
$ws.Range("E12").Value = "'7.36%"
$ws.Range("D13").Value = "'0.07972"
$ws.Range("E13").Value = "'3.82%"
$ws.Range("D14").Value = "'0.02967"
$ws.Range("E14").Value = "'1.59%"
$ws.Range("D15").Value = "'0.09910"
$ws.Range("E15").Value = "'10.33%"
$ws.Range("D16").Value = "'0.001498"
$ws.Range("E16").Value = "'-4.78%"
$ws.Range("D17").Value = "'0.04532"
$ws.Range("D18").Value = "'0.006247"
$ws.Range("E18").Value = "'-0.56%"
$ws.Range("D19").Value = "'3.489"
$ws.Range("E19").Value = "'1.10%"
$ws.Range("E20").Value = "'-0.03%"
$ws.Range("D21").Value = "'0.3339"
$ws.Range("D22").Value = "'0.1324"
$ws.Range("E22").Value = "'0.83%"
$ws.Range("D23").Value = "'4.470"
$ws.Range("E23").Value = "'9.12%"
$ws.Range("D25").Value = "'0.001214"
$ws.Range("E25").Value = "'1.59%"
$ws.Range("D26").Value = "'0.004423"
$ws.Range("E26").Value = "'6.87%"
$ws.Range("D27").Value = "'0.0001297"
$ws.Range("E27").Value = "'8.16%"
$ws.Range("D28").Value = "'0.0001739"
$ws.Range("E28").Value = "'7.56%"
$ws.Range("D40").Value = "'0.04485"
$ws.Range("E40").Value = "'5.52%"
$ws.Range("D41").Value = "'0.007195"
$ws.Range("E41").Value = "'5.14%"
$ws.Range("D42").Value = "'0.1344"
$ws.Range("E42").Value = "'8.48%"
$ws.Range("D43").Value = "'0.002324"
$ws.Range("E43").Value = "'4.32%"
$ws.Range("D44").Value = "'0.01429"
$ws.Range("E44").Value = "'10.73%"
$ws.Range("D45").Value = "'0.00006140"
$ws.Range("E45").Value = "'7.97%"
$ws.Range("D46").Value = "'1.893"
$ws.Range("E46").Value = "'-3.70%"
$ws.Range("D47").Value = "'0.01298"
$ws.Range("E47").Value = "'-13.90%"
